$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with a footnote on row 114 (empty A114, text in
# B114). A new day's worth of data needs to be appended above that
# footnote, so insert a fresh row at 114: this pushes the footnote row
# down to 115 and the new row inherits number formatting/styles from the
# row above it (row 113).
$ws.Range("A114:E114").Insert()

# Fill in the newly inserted row 114 with the new day's figures.
$ws.Cells.Item(114, 1).Value = 43969
$ws.Cells.Item(114, 2).Value = 287
$ws.Cells.Item(114, 3).Value = 37965
$ws.Cells.Item(114, 4).Value = 58
$ws.Cells.Item(114, 5).Value = 7642

# Move the active selection down to follow the newly appended data / the
# footnote row, matching where a user would land after typing this row.
$ws.Range("B115").Select() | Out-Null
